# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps for the first data row
# on the per-locale sheets, reflecting a newer report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-25 08:04:43"
$wsZhCn.Range("H2").Value = "2016-03-25 08:05:11"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-25 08:04:47"
$wsDeDe.Range("H2").Value = "2016-03-25 08:05:19"
